# Update contract_file from dashboard
# Applies the edits described by the upstream diff to the "data kontrak" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Row 26 ("030/.../Pembangunan Tambak Udang ... Tahap 1B") updates:
#    End Date, STATUS (ACTIVE -> NON ACTIVE), PROGRESS FINANCE/ACTUAL
# -----------------------------------------------------------------
$ws.Range("D26").Value = 45848
$ws.Range("D26").NumberFormat = "m/d/yyyy"
$ws.Range("I26").Value = "NON ACTIVE"
$ws.Range("M26").Value = 20
$ws.Range("N26").Value = 20

# -----------------------------------------------------------------
# 2) Insert a brand-new row at position 39 for contract
#    "031/INA02/V/2025/L-C Perjanjian Kerja Sewa Menyewa Kantor Teuku Umar"
#    This pushes the former rows 39-49 down to 40-50.
# -----------------------------------------------------------------
$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = 31
$ws.Range("B39").Value = "031/INA02/V/2025/L-C Perjanjian Kerja Sewa Menyewa Kantor Teuku Umar"
$ws.Range("C39").Value = 45778
$ws.Range("C39").NumberFormat = "m/d/yyyy"

# D39 holds the literal text "30/4/2025" (not a true date serial), right aligned
# with a date-like number format, matching the source workbook.
$ws.Range("D39").HorizontalAlignment = -4152
$ws.Range("D39").NumberFormat = "m/d/yyyy"
$ws.Range("D39").Value = "30/4/2025"

$ws.Range("F39").Value = 730
$ws.Range("G39").Value = 79
$ws.Range("I39").Value = "ACTIVE"
$ws.Range("J39").Value = "INCA & GST"
$ws.Range("M39").Value = 348

# -----------------------------------------------------------------
# 3) Row 44 (formerly the "035/.../Addendum I ... Rumah Lunyuk" row,
#    now shifted down from 43->44) is fully rewritten to describe
#    "035/.../Addendum II Perjanjian Pekerjaan Swadaya Warga".
# -----------------------------------------------------------------
$ws.Range("B44").Value = "035/INA02/V/2025/L-C Addendum II Perjanjian Pekerjaan Swadaya Warga"
$ws.Range("C44").Value = 45545
$ws.Range("C44").NumberFormat = "m/d/yyyy"
$ws.Range("D44").Value = 45940
$ws.Range("D44").NumberFormat = "m/d/yyyy"
$ws.Range("F44").Value = 425
$ws.Range("G44").Value = 311
$ws.Range("I44").Value = "ACCTIVE ADDENDUM"
$ws.Range("J44").Value = "INCA & EBS"
$ws.Range("M44").Value = 3120
$ws.Range("N44").ClearContents()

# -----------------------------------------------------------------
# 4) Row 47 (formerly "038/.../Addendum II ... Tambak Udang Sumbawa
#    Fasilitas Penunjang Swadaya Warga", shifted down from 46->47) is
#    rewritten to describe "038/.../Addendum II ... ABT dan Izin
#    Operasional". Dates/planned days/value stay the same.
# -----------------------------------------------------------------
$ws.Range("B47").Value = "038/INA02/VI/2025/L-C Addendum II Perjanjian Kerja (Kontrak) ABT dan Izin Operasional"
$ws.Range("G47").Value = 251
$ws.Range("I47").Value = "ACCTIVE ADDENDUM"
$ws.Range("J47").Value = "INCA & GEI"

# -----------------------------------------------------------------
# 5) Append three brand-new contract rows at the bottom (51-53).
# -----------------------------------------------------------------
$ws.Range("A51").Value = 43
$ws.Range("B51").Value = "042/INA02/VII/2025/L-C Addendum I Perjanjian Kerja (Kontrak) Pekerjaan Pengadaan Barang dan Jasa Serta Layanan Pendukung Operasional Kantor"
$ws.Range("B51").Font.Bold = $true
$ws.Range("C51").Value = 45778
$ws.Range("C51").NumberFormat = "m/d/yyyy"
$ws.Range("D51").Value = 46143
$ws.Range("D51").NumberFormat = "m/d/yyyy"
$ws.Range("F51").Value = 365
$ws.Range("G51").Value = 79
$ws.Range("I51").Value = "ACCTIVE ADDENDUM"
$ws.Range("J51").Value = "INCA & SAPODIA"
$ws.Range("M51").Value = 2417

$ws.Range("A52").Value = 44
$ws.Range("B52").Value = "043/INA02/VII/2025/L-C Perjanjian Kerja (Kontrak) Pekerjaan Layanan Pendukung Monitoring UKL/UPL dan PKKPRL"
$ws.Range("B52").Font.Bold = $true
$ws.Range("C52").Value = 45852
$ws.Range("C52").NumberFormat = "m/d/yyyy"
$ws.Range("D52").Value = 46217
$ws.Range("D52").NumberFormat = "m/d/yyyy"
$ws.Range("F52").Value = 365
$ws.Range("G52").Value = 5
$ws.Range("I52").Value = "ACCIVE"
$ws.Range("J52").Value = "INCA & EIJ"
$ws.Range("M52").Value = 920

$ws.Range("A53").Value = 45
$ws.Range("B53").Value = "044/INA02/VII/2025/L-C Addendum Pengadaan dan Instalasi PJU"
$ws.Range("B53").Font.Bold = $true
$ws.Range("C53").Value = 45818
$ws.Range("C53").NumberFormat = "m/d/yyyy"
$ws.Range("D53").Value = 46063
$ws.Range("D53").NumberFormat = "m/d/yyyy"
$ws.Range("F53").Value = 245
$ws.Range("G53").Value = 39
$ws.Range("I53").Value = "ACCTIVE ADDENDUM"
$ws.Range("J53").Value = "INCA & ADT"
$ws.Range("M53").Value = 197

# -----------------------------------------------------------------
# 6) Update the view: selection moves to the last new row, matching
#    where the author was working after adding the new contracts.
# -----------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B53").Select()
